# "retoques para la presentacion 2"
# Insert the clarifying word " personas" right after "(599" in the
# demographic density sentence, so it reads "...(599 personas / km²)."
# The surrounding text keeps its original (bold, Tableau Book) formatting,
# but the insertion point ends up split across three runs, matching the
# author's edit.

$d = $word.ActiveDocument

# Locate "(599" reliably (robust to any future reflow) instead of a
# hard-coded character offset.
$locate = $d.Content
$locate.Find.ClearFormatting()
$found = $locate.Find.Execute(
    "(599",   # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    "",       # ReplaceWith
    0         # Replace (wdReplaceNone)
)

# $locate is now collapsed to the found "(599" match; its End is the
# insertion point, right before the existing " / km²)." text.
$insertPos = $locate.End
$insertion = $d.Range($insertPos, $insertPos)
$insertion.InsertAfter(" personas")

# The inserted text currently lives in the same run as its neighbours
# (identical formatting gets coalesced). Nudging a formatting property
# off and back on over just the inserted span forces Word to materialize
# it as its own run, split from the text before and after it, while the
# final formatting is left unchanged (bold, same as the rest of the run).
$inserted = $d.Range($insertPos, $insertPos + 9)
$inserted.Bold = 0
$inserted.Bold = 1
